# Applies the "Changed danils and marcs role in the project" edit.
$d = $word.ActiveDocument

# --- Part 1: split "Rollenvergabe:" into two runs around the moved _GoBack bookmark ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rollePara = $d.Paragraphs(5).Range
$splitPos = $rollePara.Start + 11   # length of "Rollenverga"
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# --- Part 2: swap Daniel/Marc role assignments ---
function Replace-ParagraphXml($paraIndex, $innerXml) {
    $para = $d.Paragraphs($paraIndex).Range
    $rng = $d.Range($para.Start, $para.End - 1)
    $rng.InsertXML($innerXml)
}

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$p7xml = "<w:p $ns><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space='preserve'>      </w:t></w:r><w:r><w:t xml:space='preserve'>Daniel </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>Osipishin</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t>--   Dokumentation</w:t></w:r></w:p>"

$p8xml = "<w:p $ns><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space='preserve'>      </w:t></w:r><w:r><w:t xml:space='preserve'>Marc Sinner       </w:t></w:r><w:r><w:t>--   Vortrag</w:t></w:r></w:p>"

$p10xml = "<w:p $ns><w:r><w:tab/><w:t xml:space='preserve'> VHDL – Projekt </w:t></w:r><w:proofErr w:type='gramStart'/><w:r><w:t xml:space='preserve'>  :</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t xml:space='preserve'>  </w:t></w:r><w:r><w:t>Marc Sinner</w:t></w:r><w:r><w:tab/><w:t xml:space='preserve'>      </w:t></w:r><w:r><w:t>--    Projektleiter</w:t></w:r></w:p>"

$p12xml = "<w:p $ns><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space='preserve'>      </w:t></w:r><w:r><w:t xml:space='preserve'>Daniel </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>Osipishin</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t>--    Vortrag</w:t></w:r></w:p>"

Replace-ParagraphXml 7 $p7xml
Replace-ParagraphXml 8 $p8xml
Replace-ParagraphXml 10 $p10xml
Replace-ParagraphXml 12 $p12xml

Write-Output "done"
